{"js": "// Auto-g\u00e9n\u00e9ration des classes et des specs \u2014 apply schema/spec renames\n// and formatting fixes to the RS-RI tables.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// --- Table with \"ressource\" fields (dateTime, resourceID, RSDRId, missionID, orgID, ...) ---\nconst resourceTable = tables.items[1];\n\nconst dateTimeCell = resourceTable.getCell(1, 0);\ndateTimeCell.body.insertText(\"datetime\", Word.InsertLocation.replace);\n\nconst resourceIdCell = resourceTable.getCell(2, 0);\nresourceIdCell.body.insertText(\"resourceId\", Word.InsertLocation.replace);\n\nconst requestIdCell = resourceTable.getCell(3, 0);\nrequestIdCell.body.insertText(\"requestId\", Word.InsertLocation.replace);\n\nconst missionIdCell = resourceTable.getCell(4, 0);\nmissionIdCell.body.insertText(\"missionId\", Word.InsertLocation.replace);\n\nconst orgIdCell = resourceTable.getCell(5, 0);\norgIdCell.body.insertText(\"orgId\", Word.InsertLocation.replace);\n\nawait context.sync();\n\n// --- Table with vehicle \"status\" fields (dateTime, status, availability) ---\nconst statusTable = tables.items[3];\n\n// dateTime -> datetime\nconst statusDateTimeCell = statusTable.getCell(1, 0);\nstatusDateTimeCell.body.insertText(\"datetime\", Word.InsertLocation.replace);\n\n// dateTime row cardinality 0..1 -> 1..1\nconst statusDateTimeCardCell = statusTable.getCell(1, 3);\nstatusDateTimeCardCell.body.insertText(\"1..1\", Word.InsertLocation.replace);\n\n// \"Statuts du vecteur\" -> \"Status du vecteur\"\nconst statusLabelCell = statusTable.getCell(2, 1);\nstatusLabelCell.body.insertText(\"Status du vecteur\", Word.InsertLocation.replace);\n\n// Remove the stray space before the colon in \"(ENUM : ...\" -> \"(ENUM: ...\"\nconst statusFormatCell = statusTable.getCell(2, 2);\nstatusFormatCell.body.insertText(\n  \"string\\v(ENUM: ALERTE, PARTI, ARRIVEE SUR LES LIEUX, TRANSPORT DESTINATION, ARRIVEE DESTINATION, FIN DE MEDICALISATION , QUITTE DESTINATION, RETOUR BASE, RENTREE BASE)\",\n  Word.InsertLocation.replace\n);\n\n// status row cardinality 0..1 -> 1..1\nconst statusCardCell = statusTable.getCell(2, 3);\nstatusCardCell.body.insertText(\"1..1\", Word.InsertLocation.replace);\n\n// availability format: string + ENUM list -> boolean\nconst availabilityFormatCell = statusTable.getCell(3, 2);\navailabilityFormatCell.body.insertText(\"boolean\", Word.InsertLocation.replace);\n\n// availability description: append TRUE/FALSE/VIDE mapping on new lines\nconst availabilityDescCell = statusTable.getCell(3, 4);\navailabilityDescCell.body.insertText(\n  \"Indique si le vecteur est disponible / indisponible\\vTRUE = DISPONIBLE\\vFALSE = INDISPONIBLE\\vVIDE = INCONNU\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Auto-g\u00e9n\u00e9ration des classes et des specs \u2014 apply schema/spec renames\n# and formatting fixes to the RS-RI tables.\n\n$d = $word.ActiveDocument\n\n# --- Table 2: \"ressource\" fields (dateTime, resourceID, RSDRId, missionID, orgID, ...) ---\n$resourceTable = $d.Tables.Item(2)\n\n$resourceTable.Cell(2, 1).Range.Text = \"datetime\"\n$resourceTable.Cell(3, 1).Range.Text = \"resourceId\"\n$resourceTable.Cell(4, 1).Range.Text = \"requestId\"\n$resourceTable.Cell(5, 1).Range.Text = \"missionId\"\n$resourceTable.Cell(6, 1).Range.Text = \"orgId\"\n\n# --- Table 4: vehicle \"status\" fields (dateTime, status, availability) ---\n$statusTable = $d.Tables.Item(4)\n\n# dateTime -> datetime\n$statusTable.Cell(2, 1).Range.Text = \"datetime\"\n\n# dateTime row cardinality 0..1 -> 1..1\n$statusTable.Cell(2, 4).Range.Text = \"1..1\"\n\n# \"Statuts du vecteur\" -> \"Status du vecteur\"\n$statusTable.Cell(3, 2).Range.Text = \"Status du vecteur\"\n\n# Remove the stray space before the colon in \"(ENUM : ...\" -> \"(ENUM: ...\"\n$statusTable.Cell(3, 3).Range.Text = \"string\" + [char]11 + \"(ENUM: ALERTE, PARTI, ARRIVEE SUR LES LIEUX, TRANSPORT DESTINATION, ARRIVEE DESTINATION, FIN DE MEDICALISATION , QUITTE DESTINATION, RETOUR BASE, RENTREE BASE)\"\n\n# status row cardinality 0..1 -> 1..1\n$statusTable.Cell(3, 4).Range.Text = \"1..1\"\n\n# availability format: string + ENUM list -> boolean\n$statusTable.Cell(4, 3).Range.Text = \"boolean\"\n\n# availability description: append TRUE/FALSE/VIDE mapping on new lines\n$statusTable.Cell(4, 5).Range.Text = \"Indique si le vecteur est disponible / indisponible\" + [char]11 + \"TRUE = DISPONIBLE\" + [char]11 + \"FALSE = INDISPONIBLE\" + [char]11 + \"VIDE = INCONNU\"\n"}
